$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.615.51'
$ws.Range("E2").Value = '  +3.02%  '

$ws.Range("D3").Value = '2.545.97'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.08'
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.80'
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.555'
$ws.Range("E9").Value = '  +2.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.16'
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.38'
$ws.Range("E11").Value = '  +0.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0815'
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("E14").Value = '  +0.57%  '

$ws.Range("D15").Value = '2.944.41'
$ws.Range("E15").Value = '  +1.52%  '

$ws.Range("D16").Value = '2.546.45'
$ws.Range("E16").Value = '  +1.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.859'
$ws.Range("E17").Value = '  +1.74%  '

$ws.Range("D18").Value = '49.472.85'
$ws.Range("E18").Value = '  +3.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.06'
$ws.Range("E19").Value = '  +11.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.17'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.66'
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("E22").Value = '  -0.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '285.16'
$ws.Range("E23").Value = '  +3.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.88'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.36'

$ws.Range("E27").Value = '  -0.19%  '

$ws.Range("E28").Value = '  +6.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.143'
$ws.Range("E29").Value = '  +4.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.83'
$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.11'
$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.51'
$ws.Range("E32").Value = '  +0.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.62'
$ws.Range("E33").Value = '  +1.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.37'
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0784'
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("E37").Value = '  +2.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.69'
$ws.Range("E38").Value = '  +1.87%  '

$ws.Range("E39").Value = '  +1.16%  '

$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '121.86'
$ws.Range("E41").Value = '  -1.12%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.25'
$ws.Range("E42").Value = '  +3.14%  '

$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0311'
$ws.Range("E44").Value = '  +1.96%  '

$ws.Range("E45").Value = '  +4.88%  '

$ws.Range("D46").Value = '2.009.40'
$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("E47").Value = '  +8.43%  '

$ws.Range("E48").Value = '  +7.38%  '

$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.31'
$ws.Range("E50").Value = '  +2.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.29'
$ws.Range("E51").Value = '  +1.93%  '
